$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were added for "Ají" (Inferno, Primera/Segunda)
# dated 44900. They land at rows 92-93, pushing the previously-existing
# rows 92-114 down to 94-116 (xlShiftDown == -4121).
$ws.Rows("92:93").Insert(-4121)

# Row 92: Inferno / Primera
$ws.Cells.Item(92,1).Value = 1
$ws.Cells.Item(92,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(92,3).Value = "Arica y Parinacota"
$ws.Cells.Item(92,4).Value = 44900
$ws.Cells.Item(92,5).Value = 15
$ws.Cells.Item(92,6).Value = 100112021
$ws.Cells.Item(92,7).Value = "Ají"
$ws.Cells.Item(92,8).Value = "Inferno"
$ws.Cells.Item(92,9).Value = "Primera"
$ws.Cells.Item(92,10).Value = 300
$ws.Cells.Item(92,11).Value = 9000
$ws.Cells.Item(92,12).Value = 10000
$ws.Cells.Item(92,13).Value = 9667
$ws.Cells.Item(92,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(92,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(92,16).Value = 644
$ws.Cells.Item(92,17).Value = 15
$ws.Cells.Item(92,18).Value = "Hortaliza"

# Row 93: Inferno / Segunda
$ws.Cells.Item(93,1).Value = 1
$ws.Cells.Item(93,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(93,3).Value = "Arica y Parinacota"
$ws.Cells.Item(93,4).Value = 44900
$ws.Cells.Item(93,5).Value = 15
$ws.Cells.Item(93,6).Value = 100112021
$ws.Cells.Item(93,7).Value = "Ají"
$ws.Cells.Item(93,8).Value = "Inferno"
$ws.Cells.Item(93,9).Value = "Segunda"
$ws.Cells.Item(93,10).Value = 180
$ws.Cells.Item(93,11).Value = 6000
$ws.Cells.Item(93,12).Value = 7000
$ws.Cells.Item(93,13).Value = 6556
$ws.Cells.Item(93,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(93,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(93,16).Value = 437
$ws.Cells.Item(93,17).Value = 15
$ws.Cells.Item(93,18).Value = "Hortaliza"

Write-Output "inserted rows 92-93"
